$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(unit_sales, price) updated values for column C (unit_sales) and F (price)
$updates = @{
    4  = @(12595, 4.5)
    7  = @(13303, 4.99)
    9  = @(12668, 4.5)
    15 = @(13359, 4.99)
    16 = @(13192, 4.99)
    17 = @(17444, 4.99)
    18 = @(12422, 4.5)
    20 = @(13500, 4.99)
    21 = @(13224, 4.99)
    22 = @(13565, 4.99)
    23 = @(13134, 4.99)
    25 = @(13303, 4.99)
    26 = @(13582, 4.99)
    28 = @(13512, 4.99)
    32 = @(17185, 4.99)
    43 = @(13266, 4.99)
    46 = @(13520, 4.99)
    47 = @(13102, 4.99)
    51 = @(13314, 4.99)
    52 = @(12540, 4.5)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
}
